# LDLC Suivi Smartphones - update price history
#
# A new snapshot column is inserted right before the "nom" / "url_produit"
# columns (which sit at the end of the sheet, currently GF:GG). The new
# column:
#   - gets the timestamp "2026-02-05 16:32:45" in the header row (row 1)
#   - repeats the last known price (column GE, the previous snapshot) for
#     every product row that already had a price (rows 2-80)
#   - stays empty for the remaining rows (81-210), same as the other empty
#     snapshot columns in that range
#
# Inserting the column shifts the existing "nom" (GF) and "url_produit" (GG)
# columns one place to the right, to GG and GH respectively, and Excel keeps
# all of their data and formatting intact automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at GF, shifting "nom"/"url_produit" (and
# everything to their right) one column to the right.
$ws.Columns("GF:GF").Insert()

# New header timestamp for the inserted snapshot column.
$ws.Range("GF1").Value = "2026-02-05 16:32:45"

# Carry the previous snapshot's price (column GE) into the new column GF
# for every row that currently has a price.
$ws.Range("GF2:GF80").Value = $ws.Range("GE2:GE80").Value2

# Rows 81-210 have no price yet in any snapshot column, so GF stays blank
# there too (nothing to do - the inserted column is already empty).
